# Mac-Local-TradingModel / Temp for 2021.11.14
# Append four new daily TotalCapital rows (11/12 - 11/15/2021) below the
# existing data, re-using the "mid-series" date number format on the old
# last row (A7) and moving the "final row" date format onto the new last
# row (A11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 was previously the last row (distinct number format); now that more
# rows follow it, it takes on the same format as the other interior rows.
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New rows 8-10: interior rows.
$ws.Range("A8").Value = 44512
$ws.Range("A8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B8").Value = 55473

$ws.Range("A9").Value = 44513
$ws.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B9").Value = 55473

$ws.Range("A10").Value = 44514
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B10").Value = 55473

# New row 11: new last row, keeps the distinct "final row" number format.
$ws.Range("A11").Value = 44515
$ws.Range("A11").NumberFormat = "YYYY-MM-DD"
$ws.Range("B11").Value = 55473
